$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D28").Value = "[RLHF] Open Problems and Fundamental Limitations of Reinforcement Learning from Human Feedback"
$ws.Range("E28").Value = "https://ropiens.tistory.com/237"

$ws.Range("D32").Value = "SHAP 그래프 해석"
$ws.Range("E32").Value = "https://dodonam.tistory.com/475"

$ws.Range("D36").Value = "Diffusion-based Anomaly Detection"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/437"
